$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list data (Price + Volume(1h) columns) for rows 2-51.
# Some Price values are plain decimal-looking strings (e.g. "0.999") that Excel
# would otherwise auto-convert to numbers; force those cells to Text format first
# so the stored value stays a literal string, matching the source data feed.

$ws.Range('D2').Value = '66.306.13'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '3.559.55'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '605.15'
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.16'
$ws.Range('E6').Value = '  +1.83%  '
$ws.Range('D7').Value = '3.559.73'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -0.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.497'
$ws.Range('E9').Value = '  +1.62%  '
$ws.Range('E10').Value = '  -1.55%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.83'
$ws.Range('E11').Value = '  +0.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.410'
$ws.Range('E12').Value = '  -0.72%  '
$ws.Range('D13').Value = '4.162.54'
$ws.Range('E13').Value = '  +0.01%  '
$ws.Range('E14').Value = '  -2.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '29.21'
$ws.Range('E15').Value = '  -3.73%  '
$ws.Range('D16').Value = '3.565.27'
$ws.Range('E16').Value = '  +0.15%  '
$ws.Range('E17').Value = '  +2.03%  '
$ws.Range('D18').Value = '66.253.33'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.06'
$ws.Range('E19').Value = '  -3.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.27'
$ws.Range('E20').Value = '  +1.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.79'
$ws.Range('E21').Value = '  -0.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '420.12'
$ws.Range('E22').Value = '  -2.67%  '
$ws.Range('E23').Value = '  -1.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '77.92'
$ws.Range('E24').Value = '  -1.97%  '
$ws.Range('D25').Value = '3.699.81'
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('E27').Value = '  -1.52%  '
$ws.Range('E28').Value = '  +0.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.95'
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('E30').Value = '  -0.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('D32').Value = '3.555.79'
$ws.Range('E32').Value = '  +0.09%  '
$ws.Range('E33').Value = '  +1.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '24.64'
$ws.Range('E34').Value = '  -3.17%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.66'
$ws.Range('E36').Value = '  -2.38%  '
$ws.Range('E37').Value = '  -8.30%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.34'
$ws.Range('E38').Value = '  -4.74%  '
$ws.Range('E39').Value = '  -6.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '173.70'
$ws.Range('E40').Value = '  -1.23%  '
$ws.Range('E41').Value = '  -2.46%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.12'
$ws.Range('E42').Value = '  -1.55%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.866'
$ws.Range('E43').Value = '  -2.54%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '45.72'
$ws.Range('E44').Value = '  -0.62%  '
$ws.Range('E45').Value = '  -5.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.998'
$ws.Range('E46').Value = '  -0.14%  '
$ws.Range('E47').Value = '  -3.59%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.12'
$ws.Range('E48').Value = '  -0.32%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.01'
$ws.Range('E50').Value = '  -6.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.35'
$ws.Range('E51').Value = '  -7.24%  '
